$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 935.5106
$ws.Range("I15").Value = 935.5106
$ws.Range("K15").Value = 2806.5318
$ws.Range("M15").Value = -2637.5318
$ws.Range("H38").Value = 2705.8333
$ws.Range("I38").Value = 109
$ws.Range("J38").Value = 5302.6665
$ws.Range("K38").Value = 327
$ws.Range("L38").Value = 15907.9995
$ws.Range("M38").Value = 45
$ws.Range("N38").Value = -16651.9995
$ws.Range("H39").Value = 1123.1
$ws.Range("I39").Value = 810.2857
$ws.Range("K39").Value = 2430.8571
$ws.Range("M39").Value = -2134.8571
$ws.Range("H42").Value = 697.4286
$ws.Range("I42").Value = 57.11111
$ws.Range("J42").Value = 1850
$ws.Range("K42").Value = 171.33333
$ws.Range("L42").Value = 5550
$ws.Range("M42").Value = 58.66667000000001
$ws.Range("N42").Value = -6010
$ws.Range("H76").Value = 5553.6924
$ws.Range("I76").Value = 3858.1667
$ws.Range("J76").Value = 7007
$ws.Range("K76").Value = 3858.1667
$ws.Range("L76").Value = 7007
$ws.Range("M76").Value = -3543.1667
$ws.Range("N76").Value = -7637
$ws.Range("H79").Value = 5553.6924
$ws.Range("I79").Value = 3858.1667
$ws.Range("J79").Value = 7007
$ws.Range("K79").Value = 3858.1667
$ws.Range("L79").Value = 7007
$ws.Range("M79").Value = -2766.1667
$ws.Range("N79").Value = -9191
$ws.Range("H99").Value = 1020.35297
$ws.Range("I99").Value = 473.2857
$ws.Range("K99").Value = 1419.8571
$ws.Range("M99").Value = 78.14289999999983
$ws.Range("H135").Value = 50001890
$ws.Range("I135").Value = 52633460
$ws.Range("K135").Value = 473701140
$ws.Range("M135").Value = -473698605
$ws.Range("H137").Value = 990225.9399999999
$ws.Range("I137").Value = 1831.3334
$ws.Range("K137").Value = 5494.0002
$ws.Range("M137").Value = -2944.0002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 268.2
$ws.Range("I5").Value = 223
$ws.Range("K5").Value = 223
$ws.Range("M5").Value = -111
$ws.Range("H12").Value = 31999.666
$ws.Range("I12").Value = 22999.5
$ws.Range("J12").Value = 50000
$ws.Range("K12").Value = 22999.5
$ws.Range("L12").Value = 50000
$ws.Range("M12").Value = -22826.5
$ws.Range("N12").Value = -50346
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = ""
$ws.Range("M13").Value = ""
$ws.Range("N13").Value = 0
$ws.Range("H17").Value = 4000
$ws.Range("J17").Value = 4000
$ws.Range("L17").Value = 4000
$ws.Range("N17").Value = -4346
$ws.Range("H25").Value = 824.8333
$ws.Range("I25").Value = 829.8
$ws.Range("K25").Value = 829.8
$ws.Range("M25").Value = -427.8
$ws.Range("H32").Value = 5686652.5
$ws.Range("I32").Value = 5885876
$ws.Range("K32").Value = 5885876
$ws.Range("M32").Value = -5885589
$ws.Range("H61").Value = 14528.19
$ws.Range("I61").Value = 14054.571
$ws.Range("K61").Value = 14054.571
$ws.Range("M61").Value = -13842.571
$ws.Range("H63").Value = 7682.5
$ws.Range("I63").Value = 2299.75
$ws.Range("J63").Value = 10373.875
$ws.Range("K63").Value = 2299.75
$ws.Range("L63").Value = 10373.875
$ws.Range("M63").Value = -1613.75
$ws.Range("N63").Value = -11745.875
$ws.Range("H66").Value = 7682.5
$ws.Range("I66").Value = 2299.75
$ws.Range("J66").Value = 10373.875
$ws.Range("K66").Value = 11498.75
$ws.Range("L66").Value = 51869.375
$ws.Range("M66").Value = -8066.75
$ws.Range("N66").Value = -58733.375
$ws.Range("H132").Value = 2768.4902
$ws.Range("I132").Value = 2358.9092
$ws.Range("J132").Value = 5343
$ws.Range("K132").Value = 7076.7276
$ws.Range("L132").Value = 16029
$ws.Range("M132").Value = -4546.7276
$ws.Range("N132").Value = -21089
$ws.Range("H136").Value = 14528.19
$ws.Range("I136").Value = 14054.571
$ws.Range("K136").Value = 42163.713
$ws.Range("M136").Value = -39613.713

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 268.2
$ws.Range("I4").Value = 223
$ws.Range("K4").Value = 223
$ws.Range("M4").Value = -108
$ws.Range("H5").Value = 4847.1665
$ws.Range("I5").Value = 876
$ws.Range("J5").Value = 8818.333000000001
$ws.Range("K5").Value = 876
$ws.Range("L5").Value = 8818.333000000001
$ws.Range("M5").Value = -763
$ws.Range("N5").Value = -9044.333000000001
$ws.Range("H7").Value = 1500
$ws.Range("J7").Value = 1500
$ws.Range("L7").Value = 1500
$ws.Range("N7").Value = -1726
$ws.Range("H11").Value = 3966.6667
$ws.Range("I11").Value = 1900
$ws.Range("K11").Value = 1900
$ws.Range("M11").Value = -1760
$ws.Range("H12").Value = 2820
$ws.Range("I12").Value = 180
$ws.Range("K12").Value = 180
$ws.Range("M12").Value = -12
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = ""
$ws.Range("M17").Value = ""
$ws.Range("N17").Value = 0
$ws.Range("H19").Value = 900
$ws.Range("I19").Value = 900
$ws.Range("K19").Value = 900
$ws.Range("M19").Value = -727
$ws.Range("H23").Value = 1499
$ws.Range("J23").Value = 1499
$ws.Range("L23").Value = 1499
$ws.Range("N23").Value = -2065
$ws.Range("H25").Value = 4514
$ws.Range("I25").Value = 4514
$ws.Range("K25").Value = 4514
$ws.Range("M25").Value = -4279
$ws.Range("H134").Value = 2667.2927
$ws.Range("I134").Value = 2186.0715
$ws.Range("J134").Value = 5474.4165
$ws.Range("K134").Value = 6558.2145
$ws.Range("L134").Value = 16423.2495
$ws.Range("M134").Value = -4023.2145
$ws.Range("N134").Value = -21493.2495

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5150.4653
$ws.Range("I31").Value = 1858.0968
$ws.Range("J31").Value = 8930.593000000001
$ws.Range("K31").Value = 1858.0968
$ws.Range("L31").Value = 8930.593000000001
$ws.Range("M31").Value = -1563.0968
$ws.Range("N31").Value = -9520.593000000001
$ws.Range("H34").Value = 5150.4653
$ws.Range("I34").Value = 1858.0968
$ws.Range("J34").Value = 8930.593000000001
$ws.Range("K34").Value = 1858.0968
$ws.Range("L34").Value = 8930.593000000001
$ws.Range("M34").Value = -1656.0968
$ws.Range("N34").Value = -9334.593000000001
$ws.Range("H58").Value = 2691.6365
$ws.Range("J58").Value = 4249.75
$ws.Range("L58").Value = 4249.75
$ws.Range("N58").Value = -4655.75
$ws.Range("H68").Value = 86515.625
$ws.Range("J68").Value = 86515.625
$ws.Range("L68").Value = 86515.625
$ws.Range("N68").Value = -88013.625
$ws.Range("H71").Value = 86515.625
$ws.Range("J71").Value = 86515.625
$ws.Range("L71").Value = 259546.875
$ws.Range("N71").Value = -267034.875
$ws.Range("H74").Value = 46535.312
$ws.Range("J74").Value = 46665.332
$ws.Range("L74").Value = 46665.332
$ws.Range("N74").Value = -48413.332
$ws.Range("H77").Value = 46535.312
$ws.Range("J77").Value = 46665.332
$ws.Range("L77").Value = 139995.996
$ws.Range("N77").Value = -148731.996
$ws.Range("H132").Value = 11906945
$ws.Range("I132").Value = 2108.6667
$ws.Range("K132").Value = 6326.000100000001
$ws.Range("M132").Value = -3796.000100000001
$ws.Range("H136").Value = 2691.6365
$ws.Range("J136").Value = 4249.75
$ws.Range("L136").Value = 12749.25
$ws.Range("N136").Value = -17849.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 451.8889
$ws.Range("I33").Value = 525.25
$ws.Range("J33").Value = 393.2
$ws.Range("K33").Value = 3151.5
$ws.Range("L33").Value = 2359.2
$ws.Range("M33").Value = -2868.5
$ws.Range("N33").Value = -2925.2
$ws.Range("H92").Value = 254
$ws.Range("I92").Value = 326.6
$ws.Range("J92").Value = 181.4
$ws.Range("K92").Value = 979.8000000000001
$ws.Range("L92").Value = 544.2
$ws.Range("M92").Value = 268.1999999999999
$ws.Range("N92").Value = -3040.2
$ws.Range("H102").Value = 7738.222
$ws.Range("I102").Value = 3846
$ws.Range("J102").Value = 9684.333000000001
$ws.Range("K102").Value = 11538
$ws.Range("L102").Value = 29052.999
$ws.Range("M102").Value = -9104
$ws.Range("N102").Value = -33920.999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 9398
$ws.Range("I122").Value = 11897.143
$ws.Range("J122").Value = 3566.6667
$ws.Range("K122").Value = 35691.429
$ws.Range("L122").Value = 10700.0001
$ws.Range("M122").Value = -33241.429
$ws.Range("N122").Value = -15600.0001
$ws.Range("H132").Value = 288888.5
$ws.Range("I132").Value = 483222.03
$ws.Range("J132").Value = 3103.8823
$ws.Range("K132").Value = 1449666.09
$ws.Range("L132").Value = 9311.6469
$ws.Range("M132").Value = -1447136.09
$ws.Range("N132").Value = -14371.6469

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 7334.2856
$ws.Range("I100").Value = 2684.8462
$ws.Range("K100").Value = 2684.8462
$ws.Range("M100").Value = -2143.8462
$ws.Range("H132").Value = 825969.2
$ws.Range("I132").Value = 1082744
$ws.Range("J132").Value = 4289.9
$ws.Range("K132").Value = 3248232
$ws.Range("L132").Value = 12869.7
$ws.Range("M132").Value = -3245702
$ws.Range("N132").Value = -17929.7

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 482516
$ws.Range("I132").Value = 753620.2
$ws.Range("J132").Value = 2870.1155
$ws.Range("K132").Value = 2260860.6
$ws.Range("L132").Value = 8610.3465
$ws.Range("M132").Value = -2258330.6
$ws.Range("N132").Value = -13670.3465
